$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tier 1_obs")

for ($r = 6; $r -le 34; $r++) {
    $iCell = $ws.Cells.Item($r, 9)        # column I
    $kCell = $ws.Cells.Item($r, 11)       # column K

    $iVal = $iCell.Value()
    if ($iVal -ne $null -and $iVal -ne "") {
        $ws.Cells.Item($r, 10).Value = $iVal   # column J = column I
    }

    $kVal = $kCell.Value()
    if ($kVal -ne $null -and $kVal -ne "" -and $kVal -notlike "NCBITaxon:*") {
        $kCell.Value = "NCBITaxon:" + $kVal
    }
}
